# Update scripts with new TPM (transcripts-per-million) derived NATMI
# ligand-receptor output for the Oxt -> Avpr2 pair. Target cluster moves
# from "Resolving-Mac" to "Inflammatory-Mac" for every row, the third
# sending cluster becomes "MuSCs" (was "Inflammatory-Mac"), and all of
# the recomputed detection / expression / specificity metrics are
# refreshed to match the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2955753333333333
$ws.Range("H2").Value = 0.8867259999999999
$ws.Range("I2").Value = 0.239018529794766
$ws.Range("J2").Value = 0.2584571780171812
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09302566666666667
$ws.Range("N2").Value = 0.279077
$ws.Range("Q2").Value = 0.02749609243355556
$ws.Range("R2").Value = 0.247464831902
$ws.Range("S2").Value = 0.239018529794766
$ws.Range("T2").Value = 0.2584571780171812
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("I3").Value = 0.2089742936599006
$ws.Range("J3").Value = 0.2259695357671569
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09302566666666667
$ws.Range("N3").Value = 0.279077
$ws.Range("Q3").Value = 0.02403987883133333
$ws.Range("R3").Value = 0.216358909482
$ws.Range("S3").Value = 0.2089742936599006
$ws.Range("T3").Value = 0.2259695357671569
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.27902
$ws.Range("H4").Value = 0.55804
$ws.Range("I4").Value = 0.2256309734348681
$ws.Range("J4").Value = 0.1626539016795581
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09302566666666667
$ws.Range("N4").Value = 0.279077
$ws.Range("Q4").Value = 0.02595602151333333
$ws.Range("R4").Value = 0.15573612908
$ws.Range("S4").Value = 0.2256309734348681
$ws.Range("T4").Value = 0.1626539016795581
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4036036666666667
$ws.Range("H5").Value = 1.210811
$ws.Range("I5").Value = 0.3263762031104653
$ws.Range("J5").Value = 0.3529193845361038
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09302566666666667
$ws.Range("N5").Value = 0.279077
$ws.Range("Q5").Value = 0.03754550016077778
$ws.Range("R5").Value = 0.337909501447
$ws.Range("S5").Value = 0.3263762031104653
$ws.Range("T5").Value = 0.3529193845361038
